$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.266.05'
$ws.Range("E2").Value = '  +3.73%  '

$ws.Range("D3").Value = '1.809.29'
$ws.Range("E3").Value = '  +1.46%  '

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = '  -0.36%  '

$ws.Range("D5").Value = "'339.46"
$ws.Range("E5").Value = '  +1.15%  '

$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = '  -0.44%  '

$ws.Range("D7").Value = "'0.3936"
$ws.Range("E7").Value = '  +4.19%  '

$ws.Range("D8").Value = "'0.3495"
$ws.Range("E8").Value = '  +2.50%  '

$ws.Range("D9").Value = "'48.19"
$ws.Range("E9").Value = '  -0.15%  '

$ws.Range("E10").Value = '  -1.10%  '

$ws.Range("D11").Value = "'0.07546"
$ws.Range("E11").Value = '  +1.64%  '

$ws.Range("D12").Value = "'0.9991"
$ws.Range("E12").Value = '  -0.34%  '

$ws.Range("D13").Value = "'22.08"
$ws.Range("E13").Value = '  +2.39%  '

$ws.Range("D14").Value = "'6.518"
$ws.Range("E14").Value = '  +1.93%  '

$ws.Range("D15").Value = '1.812.92'
$ws.Range("E15").Value = '  +1.90%  '

$ws.Range("D16").Value = "'7.163"
$ws.Range("E16").Value = '  +1.98%  '

$ws.Range("D17").Value = "'0.00001105"
$ws.Range("E17").Value = '  +1.51%  '

$ws.Range("D18").Value = "'0.06719"
$ws.Range("E18").Value = '  +0.70%  '

$ws.Range("D19").Value = "'85.32"
$ws.Range("E19").Value = '  +1.51%  '

$ws.Range("E20").Value = '  -0.40%  '

$ws.Range("D21").Value = "'17.74"
$ws.Range("E21").Value = '  +3.01%  '

$ws.Range("D22").Value = "'6.563"
$ws.Range("E22").Value = '  +0.64%  '

$ws.Range("D23").Value = '28.238.78'

$ws.Range("E24").Value = '  +0.38%  '

$ws.Range("D25").Value = "'2.401"

$ws.Range("D26").Value = "'21.40"
$ws.Range("E26").Value = '  +1.51%  '

$ws.Range("E27").Value = '  -1.47%  '

$ws.Range("D28").Value = "'2.523"
$ws.Range("E28").Value = '  +0.16%  '

$ws.Range("D29").Value = "'154.99"
$ws.Range("E29").Value = '  +1.75%  '

$ws.Range("D30").Value = '2.016.92'
$ws.Range("E30").Value = '  +1.71%  '

$ws.Range("D31").Value = "'136.09"
$ws.Range("E31").Value = '  +2.68%  '

$ws.Range("D32").Value = "'6.267"
$ws.Range("E32").Value = '  +4.66%  '

$ws.Range("D33").Value = "'4.017"
$ws.Range("E33").Value = '  -1.55%  '

$ws.Range("D34").Value = "'0.08847"
$ws.Range("E34").Value = '  +3.09%  '

$ws.Range("D35").Value = "'13.19"
$ws.Range("E35").Value = '  +1.36%  '

$ws.Range("D36").Value = "'0.02450"
$ws.Range("E36").Value = '  +4.88%  '

$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").Value = "'0.06572"
$ws.Range("E37").Value = '  +4.05%  '

$ws.Range("D38").Value = "'5.468"
$ws.Range("E38").Value = '  +1.40%  '

$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").Value = "'0.6931"
$ws.Range("E39").Value = '  +1.89%  '

$ws.Range("D40").Value = "'1.615"
$ws.Range("E40").Value = '  -1.83%  '

$ws.Range("D41").Value = "'0.2229"
$ws.Range("E41").Value = '  +2.54%  '

$ws.Range("D42").Value = "'1.264"
$ws.Range("E42").Value = '  +1.51%  '

$ws.Range("D43").Value = "'8.566"
$ws.Range("E43").Value = '  -2.21%  '

$ws.Range("D44").Value = "'14.62"
$ws.Range("E44").Value = '  +1.65%  '

$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").Value = "'1.000"
$ws.Range("E45").Value = '  -0.44%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = "'0.6418"
$ws.Range("E46").Value = '  +1.18%  '

$ws.Range("D47").Value = "'3.881"
$ws.Range("E47").Value = '  +1.01%  '

$ws.Range("D48").Value = "'2.156"
$ws.Range("E48").Value = '  +2.02%  '

$ws.Range("D49").Value = "'131.47"
$ws.Range("E49").Value = '  +2.26%  '

$ws.Range("D50").Value = "'0.07216"
$ws.Range("E50").Value = '  +0.79%  '

$ws.Range("D51").Value = "'80.27"
$ws.Range("E51").Value = '  +1.70%  '
